$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.585.33"
$ws.Range("E2").Value = "  +3.76%  "

$ws.Range("D3").Value = "1.796.77"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'313.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").Value = "'0.5285"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.42%  "

$ws.Range("D8").Value = "'0.3769"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.13%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'42.65"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.51%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.07520"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.09%  "

$ws.Range("D11").Value = "'1.117"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.29%  "

$ws.Range("B12").Value = "BinanceUSD"
$ws.Range("C12").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D12").Value = "'1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'21.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.04%  "

$ws.Range("D14").Value = "'7.505"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.93%  "

$ws.Range("D15").Value = "'6.183"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "1.796.17"
$ws.Range("E16").Value = "  +0.45%  "

$ws.Range("D17").Value = "'90.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.63%  "

$ws.Range("D18").Value = "'0.00001065"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").Value = "'0.06465"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("D21").Value = "'17.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.86%  "

$ws.Range("E22").Value = "  -0.27%  "

$ws.Range("D23").Value = "28.592.03"
$ws.Range("E23").Value = "  +3.59%  "

$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").Value = "'2.092"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.17%  "

$ws.Range("D26").Value = "'160.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.54%  "

$ws.Range("D27").Value = "'20.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.19%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.367"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.07%  "

$ws.Range("B29").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C29").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D29").Value = "2.001.05"
$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("D30").Value = "'123.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.64%  "

$ws.Range("D31").Value = "'1.113"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.67%  "

$ws.Range("D32").Value = "'0.1029"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "

$ws.Range("D33").Value = "'5.692"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.23%  "

$ws.Range("D34").Value = "'3.678"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.95%  "

$ws.Range("D35").Value = "'0.2278"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.39%  "

$ws.Range("D36").Value = "'0.06510"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.90%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "'8.899"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.66%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.02314"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.45%  "

$ws.Range("D39").Value = "'5.057"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.36%  "

$ws.Range("D40").Value = "'11.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.39%  "

$ws.Range("D41").Value = "'0.6277"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.73%  "

$ws.Range("D42").Value = "'1.212"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.69%  "

$ws.Range("D43").Value = "'1.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.25%  "

$ws.Range("D44").Value = "'1.396"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.43%  "

$ws.Range("D45").Value = "'13.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.49%  "

$ws.Range("D46").Value = "'0.5908"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.73%  "

$ws.Range("D47").Value = "'3.667"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.87%  "

$ws.Range("D48").Value = "'126.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.25%  "

$ws.Range("D49").Value = "'1.973"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.93%  "

$ws.Range("D50").Value = "'1.161"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.50%  "

$ws.Range("D51").Value = "'0.06922"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.60%  "
